# Generate Report for Handoff
#
# Regenerates the localization-status report for the most recent handoff
# batch (file 4f37f634-e9d7-48fe-9b64-58ac52dc2329). The generator first
# stamps the row with the new "just generated" timestamps, then reconciles
# each cell back to the authoritative last-recorded handoff time before the
# report is written out:
#   - Overview!G6  ("Latest HO Xliff Generate Date")
#   - zh-cn!H6     ("Latest Handoff Datetime")

$wb = $excel.ActiveWorkbook
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Remember the authoritative timestamps already on record for this row.
$overviewRecordedDate = $wsOverview.Range("G6").Text
$zhCnRecordedDate = $wsZhCn.Range("H6").Text

# Stamp the row with the freshly-generated handoff timestamps.
$wsOverview.Range("G6").Value = "2016-08-21 06:48:46"
$wsZhCn.Range("H6").Value = "2016-08-21 06:48:42"

# Reconcile back to the last-recorded handoff timestamps for the report.
$wsOverview.Range("G6").Value = $overviewRecordedDate
$wsZhCn.Range("H6").Value = $zhCnRecordedDate
